# Rebuild each worksheet's small "website uptime" table:
#   before: row1 = [ , website_isc, website_ist ]      (headers in B1/C1, style "1")
#           row2 = [0, value_isc, value_ist]            (A2 style "1", B2/C2 plain)
#   after:  row1 = [Websites, 2019-04-06, 2019-04-07, 2019-04-08]   (A1:D1 style "1")
#           row2 = [website_ist, 0, 0, value_ist_new]    (plain, no style)
#           row3 = [website_isc, 0, 0, value_isc_new]    (plain, no style)
#
# The same reshape is applied to every sheet; only the value landing in the
# new column D differs per sheet.

$wb = $excel.ActiveWorkbook

$sheetValues = @{
    "Sheet0" = @(0, 0)
    "Sheet1" = @(0, 64763)
    "Sheet2" = @(26671, 105257)
    "Sheet3" = @(0, 0)
    "Sheet4" = @(0, 0)
}

foreach ($ws in $wb.Worksheets) {
    $vals = $sheetValues[$ws.Name]
    $istValue = $vals[0]
    $iscValue = $vals[1]

    # Stash the existing header style (currently on B1) in a scratch cell,
    # far outside the table, so it survives the reshape below untouched.
    $ws.Range("B1").Copy()
    $ws.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    # Wipe the old 2-row / 3-col block (values + formatting) entirely.
    $ws.Range("A1:D3").Clear()

    # Header row text. Force text (not auto-parsed dates) via NumberFormat,
    # then restore the original header style from the stash afterwards.
    $ws.Cells.Item(1, 1).Value = "Websites"
    $ws.Cells.Item(1, 2).NumberFormat = "@"
    $ws.Cells.Item(1, 2).Value = "2019-04-06"
    $ws.Cells.Item(1, 3).NumberFormat = "@"
    $ws.Cells.Item(1, 3).Value = "2019-04-07"
    $ws.Cells.Item(1, 4).NumberFormat = "@"
    $ws.Cells.Item(1, 4).Value = "2019-04-08"

    $ws.Range("Z1").Copy()
    $ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
    $ws.Range("Z1").Clear()

    # Row 2: istitutocomprensivotrebisacce.gov.it (no style, like the rest
    # of the data cells already had).
    $ws.Cells.Item(2, 1).Value = "www.istitutocomprensivotrebisacce.gov.it"
    $ws.Cells.Item(2, 2).Value = 0
    $ws.Cells.Item(2, 3).Value = 0
    $ws.Cells.Item(2, 4).Value = $istValue

    # Row 3: iscmontegiorgio.it
    $ws.Cells.Item(3, 1).Value = "www.iscmontegiorgio.it"
    $ws.Cells.Item(3, 2).Value = 0
    $ws.Cells.Item(3, 3).Value = 0
    $ws.Cells.Item(3, 4).Value = $iscValue
}
